{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n// copyright/footer line right after it, and the blank paragraph that\n// separates them from the preceding \"Requisitos\" text, while leaving the\n// rest of the document (including the final blank / page-break paragraphs)\n// untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet footerStart = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    footerStart = i;\n    break;\n  }\n}\n\nif (footerStart > -1) {\n  // The blank paragraph immediately preceding the footer block is removed\n  // together with the two text paragraphs that follow it.\n  const toDelete = [];\n  if (footerStart - 1 >= 0 && items[footerStart - 1].text === \"\") {\n    toDelete.push(items[footerStart - 1]);\n  }\n  toDelete.push(items[footerStart]);\n  if (footerStart + 1 < items.length && items[footerStart + 1].text === targetTexts[1]) {\n    toDelete.push(items[footerStart + 1]);\n  }\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n# copyright/footer line right after it, and the blank paragraph that\n# separates them from the preceding \"Requisitos\" text, while leaving the\n# rest of the document (including the final blank / page-break paragraphs)\n# untouched.\n$d = $word.ActiveDocument\n\n$footerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightMarker = \"Contact: luizeleno@usp.br\"\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $footerText) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    # Delete the copyright paragraph right after the footer line, if present.\n    if ($targetIndex + 1 -le $d.Paragraphs.Count) {\n        $afterText = $d.Paragraphs.Item($targetIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($afterText -like \"*$copyrightMarker*\") {\n            $d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n        }\n    }\n\n    # Delete the footer line itself.\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n\n    # Delete the blank separator paragraph right before it, if present.\n    if ($targetIndex - 1 -ge 1) {\n        $beforeText = $d.Paragraphs.Item($targetIndex - 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($beforeText -eq \"\") {\n            $d.Paragraphs.Item($targetIndex - 1).Range.Delete()\n        }\n    }\n}\n"}
